# Append a new plate2 thermal-curve data block (date 20250627, temperature 26)
# to Sheet1, following the existing 29-row-per-plate pattern:
#   wells A01-A12, B01-B12 => type "sample"
#   wells C01-C05          => type "blank"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$date = 20250627
$plate = "plate2"
$temperature = 26

$wells = @(
    "A01","A02","A03","A04","A05","A06","A07","A08","A09","A10","A11","A12",
    "B01","B02","B03","B04","B05","B06","B07","B08","B09","B10","B11","B12",
    "C01","C02","C03","C04","C05"
)

$startRow = 118
$row = $startRow

foreach ($well in $wells) {
    if ($well.StartsWith("C")) {
        $type = "blank"
    } else {
        $type = "sample"
    }

    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $plate
    $ws.Cells.Item($row, 3).Value = $temperature
    $ws.Cells.Item($row, 4).Value = $well
    $ws.Cells.Item($row, 5).Value = $type

    $row = $row + 1
}

$lastRow = $row - 1

# Mirror the final selection recorded in the workbook after the edit.
$ws.Range("C$startRow`:C$lastRow").Select()
